$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9019824266433716
$ws.Range("B1").Value = 1.270700573921204
$ws.Range("C1").Value = 2.211973667144775
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 1.69395923614502
